# update eval data of LLM
# - refresh the cached eval-result JSON in column D (drop the stale
#   "completionTokens" field, keep "fcCount")
# - leave the current view scrolled/selected over the results columns

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$newPayload = '{"fcCount":1}'

$lastRow = $ws.Cells(1, 1).End(-4121).Row   # -4121 = xlDown
if ($lastRow -lt 2) { $lastRow = 51 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)   # column D
    if ($cell.Value2 -ne $null -and $cell.Value2 -ne "") {
        $cell.Value = $newPayload
    }
}

# Scroll/select like the author left the workbook: columns C.. in view,
# E2:L52 highlighted as the active selection.
$ws.Range("E2:L52").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
